$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Force column A's new cell to be stored as text (it would otherwise be
# auto-recognized as a date literal), then restore the default/"Normal"
# style so no extra cell formatting is introduced.
$ws.Range("A68").NumberFormat = "@"
$ws.Range("A68").Value = "01-07-2021"
$ws.Range("A68").Style = "Normal"

$ws.Range("B68").Value = -3.3
